# Add a new inventory row (row 69) to the bottom of the sheet, mirroring the
# layout/format of the preceding rows (66-68): Host, Grupo, IP, Usuario,
# Descripción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the row above so the new row picks up the same cell styles
# (A/C/E use the "Aptos Narrow" style, B/D stay default) instead of minting
# brand-new style entries.
$ws.Range("A66:E66").Copy($ws.Range("A69:E69"))

# Overwrite with the new host's data.
$ws.Range("A69").Value = "10.181.11.192"
$ws.Range("B69").Value = "CRM"
$ws.Range("C69").Value = "10.181.11.192"
$ws.Range("D69").Value = "tux"
$ws.Range("E69").Value = "crmtestmicros"

# Match the author's final selection after entering the new row.
$null = $ws.Range("E69").Select()
